$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update dataset links/dates per latest nomis/ILR releases and reload

# Row 5: FE and skills achievements and participation (ILR) - update link to 2022-23 and latest period
$ws.Range("B5").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/further-education-and-skills/2022-23'>Individualised Learner Record</a>"
$ws.Range("C5").Value = "Aug 2021 – Jul 2022 (31/03/23)"

# Row 2: Employment volumes (APS)
$ws.Range("B2").Value = "<a href='https://www.nomisweb.co.uk/livelinks/16244.xlsx'>Annual Population Survey</a>"

# Row 4: Employment by industry (APS)
$ws.Range("B4").Value = "<a href='https://www.nomisweb.co.uk/livelinks/16243.xlsx'>Annual Population Survey</a>"

# Row 8: Enterprises by employment size band (ONS UK Business Counts)
$ws.Range("B8").Value = "<a href='https://www.nomisweb.co.uk/livelinks/16246.xlsx'>ONS UK Business Counts</a>"

# Row 3: Employment by occupation (APS) - next period now TBC
$ws.Range("D3").Value = "TBC"

# Row 2 & 4: refreshed APS periods
$ws.Range("D2").Value = "Apr 2022 - Mar 2023 (11/06/23)"
$ws.Range("C2").Value = "Jan 2022 - Dec 2022 (19/04/23)"
$ws.Range("C4").Value = "Jan 2022 - Dec 2022 (19/04/23)"
$ws.Range("D4").Value = "Apr 2022 - Mar 2023 (11/06/23)"

# Row 7: Highest qualification level by age and gender (APS) - next period now TBC
$ws.Range("D7").Value = "TBC"

# Remove the date-format style applied to C2/C4 so it matches default (no longer a distinct date style)
$ws.Range("C2").Style = "Normal"
$ws.Range("C4").Style = "Normal"

# Update sheet view: remove frozen/topLeftCell positioning, set new selection
$ws.Range("C5").Select()
